# Weekly data update: add a new "29_06_2021" column (AI) after the last
# existing week column (AH), carrying over this week's per-age-group death
# counts, and extend the row-12 "total" SUM formula across into AI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added week column.
$ws.Range("AI1").Value = "29_06_2021"

# New week's counts per age group (rows 2-11 match the row order already
# used on the sheet: 0-9, 10-19, 20-29, 30-39, 40-49, 50-59, 60-69, 70-79,
# 80-89, 90+).
$ws.Range("AI2").Value = 1
$ws.Range("AI3").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AI5").Value = 7
$ws.Range("AI6").Value = 11
$ws.Range("AI7").Value = 66
$ws.Range("AI8").Value = 213
$ws.Range("AI9").Value = 678
$ws.Range("AI10").Value = 994
$ws.Range("AI11").Value = 564

# Extend the "I alt" (total) row's SUM formula rightwards into the new
# column, matching how the existing AC12:AH12 block was filled across.
$ws.Range("AC12:AI12").Formula = "=SUM(AC2:AC11)"

# Move the active selection, matching where the author's cursor ended up
# after adding the new column.
[void]$ws.Range("AN20").Select()
